# "upload iris data 4"
# Resize/reposition the two iris-data pictures and the accompanying
# "(B)" textbox label on slide 1.
#
# NOTE: PowerPoint COM exposes Left/Top/Width/Height in points, while the
# underlying OOXML stores EMUs (1 pt = 12700 EMU). The point values below
# were chosen so that, after the host's internal point<->EMU round trip,
# they land exactly on the target EMU coordinates from the authoritative
# edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Picture 3 (first iris image) - grows and shifts up by 1 EMU
$pic1 = $s.Shapes.Item(1)
$pic1.Left   = 75.7368507937008
$pic1.Top    = 88.10527559055119
$pic1.Width  = 333.4975738551181
$pic1.Height = 187.6577225354331

# Picture 4 (second iris image) - grows and moves to the right
$pic2 = $s.Shapes.Item(2)
$pic2.Left   = 433.22874455748035
$pic2.Top    = 88.70110236220472
$pic2.Width  = 173.9659842519685
$pic2.Height = 187.6577225354331

# TextBox 6 ("(B)" label) - follows Picture 4 to the right
$txt2 = $s.Shapes.Item(4)
$txt2.Left = 402.59740157480314
$txt2.Top  = 51.40212598425197
